$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 46 / 47: coin identity swap (ONDO <-> Stellar) ---
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.105"
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.02"

# --- Price (column D) updates ---
# Cells whose new text would otherwise be reinterpreted as a shorter number
# (trailing/implied zero lost) are pinned to Text format first.
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"

$ws.Range("D2").Value = "61.006.39"
$ws.Range("D3").Value = "3.255.82"
$ws.Range("D5").Value = "545.23"
$ws.Range("D6").Value = "148.15"
$ws.Range("D7").Value = "1.00"
$ws.Range("D8").Value = "0.523"
$ws.Range("D11").Value = "0.433"
$ws.Range("D12").Value = "3.818.94"
$ws.Range("D14").Value = "26.48"
$ws.Range("D16").Value = "61.046.82"
$ws.Range("D17").Value = "3.272.74"
$ws.Range("D18").Value = "6.34"
$ws.Range("D19").Value = "13.45"
$ws.Range("D20").Value = "8.43"
$ws.Range("D21").Value = "377.75"
$ws.Range("D22").Value = "1.00"
$ws.Range("D23").Value = "0.532"
$ws.Range("D24").Value = "70.11"
$ws.Range("D26").Value = "8.63"
$ws.Range("D28").Value = "0.0₃0923"
$ws.Range("D30").Value = "22.62"
$ws.Range("D31").Value = "6.20"
$ws.Range("D32").Value = "5.42"
$ws.Range("D33").Value = "1.25"
$ws.Range("D35").Value = "159.43"
$ws.Range("D36").Value = "1.45"
$ws.Range("D37").Value = "26.43"
$ws.Range("D38").Value = "2.791.68"
$ws.Range("D40").Value = "1.73"
$ws.Range("D41").Value = "0.0313"
$ws.Range("D42").Value = "4.28"
$ws.Range("D43").Value = "40.11"
$ws.Range("D44").Value = "0.733"
$ws.Range("D45").Value = "3.301.08"
$ws.Range("D51").Value = "278.05"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +3.08%  "
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("E6").Value = "  +4.45%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("E10").Value = "  +3.09%  "
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("E21").Value = "  +0.43%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("E25").Value = "  +2.33%  "
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +5.94%  "
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  +2.71%  "
$ws.Range("E32").Value = "  +4.03%  "
$ws.Range("E33").Value = "  +7.22%  "
$ws.Range("E34").Value = "  +4.77%  "
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +7.89%  "
$ws.Range("E37").Value = "  +4.20%  "
$ws.Range("E38").Value = "  +3.20%  "
$ws.Range("E39").Value = "  +1.19%  "
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("E41").Value = "  +6.38%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  +2.54%  "
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("E48").Value = "  +6.81%  "
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("E50").Value = "  +6.27%  "
$ws.Range("E51").Value = "  +8.35%  "

Write-Output "applied cryptos update"
